$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at A, shifting existing columns (A..AA) to (B..AB)
$ws.Columns("A:A").Insert()

# Set the column width for the new column A to match the target layout.
# (ColumnWidth is character-width; the engine adds ~5/6 char of padding when
# serializing the raw OOXML <col width>, so back the padding out here so the
# saved width attribute comes out to exactly 23.)
$ws.Columns("A:A").ColumnWidth = 22.166666666666668

# Match the cell formatting (style) of the new column A to the adjacent
# (now shifted) column B, which holds the same header/data row styles.
$ws.Range("B1:B7").Copy()
$ws.Range("A1:A7").PasteSpecial(-4122)  # xlPasteFormats

# Header for the new index column
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"

# Populate the per-row index values (these are the new "DO NOT MODIFY" ids)
$ws.Range("A2").Value = 108
$ws.Range("A3").Value = 1733
$ws.Range("A4").Value = 1734
$ws.Range("A5").Value = 1735
$ws.Range("A6").Value = 1736
$ws.Range("A7").Value = 2385

Write-Host "Done"
